$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1958333333333333
$ws.Range("C2").Value = 0.55
$ws.Range("J2").Value = 0.004166666666666667
$ws.Range("P2").Value = 0.1416666666666667
$ws.Range("S2").Value = 0.1083333333333333
$ws.Range("B3").Value = 0.007575757575757576
$ws.Range("C3").Value = 0.03787878787878788
$ws.Range("J3").Value = 0.0303030303030303
$ws.Range("P3").Value = 0.7348484848484849
$ws.Range("S3").Value = 0.1893939393939394
$ws.Range("P4").Value = 0.7058823529411765
$ws.Range("S4").Value = 0.2941176470588235
$ws.Range("B6").Value = 0.06153846153846154
$ws.Range("D6").Value = 0.03076923076923077
$ws.Range("F6").Value = 0.03589743589743589
$ws.Range("J6").Value = 0.2051282051282051
$ws.Range("O6").Value = 0.03076923076923077
$ws.Range("Q6").Value = 0.1794871794871795
$ws.Range("R6").Value = 0.04102564102564103
$ws.Range("S6").Value = 0.4153846153846154
$ws.Range("B7").Value = 0.1317829457364341
$ws.Range("D7").Value = 0.02325581395348837
$ws.Range("F7").Value = 0.02325581395348837
$ws.Range("J7").Value = 0.1395348837209302
$ws.Range("O7").Value = 0.03875968992248062
$ws.Range("Q7").Value = 0.2403100775193799
$ws.Range("R7").Value = 0.08527131782945736
$ws.Range("S7").Value = 0.3178294573643411
$ws.Range("B8").Value = 0.1181102362204724
$ws.Range("D8").Value = 0.005249343832020997
$ws.Range("F8").Value = 0.06561679790026247
$ws.Range("J8").Value = 0.07874015748031496
$ws.Range("O8").Value = 0.03149606299212598
$ws.Range("Q8").Value = 0.1916010498687664
$ws.Range("R8").Value = 0.05774278215223097
$ws.Range("S8").Value = 0.4514435695538058
$ws.Range("B9").Value = 0.09677419354838709
$ws.Range("D9").Value = 0.01209677419354839
$ws.Range("F9").Value = 0.06451612903225806
$ws.Range("J9").Value = 0.07661290322580645
$ws.Range("O9").Value = 0.0564516129032258
$ws.Range("Q9").Value = 0.157258064516129
$ws.Range("R9").Value = 0.04838709677419355
$ws.Range("S9").Value = 0.4879032258064516
$ws.Range("B10").Value = 0.0851063829787234
$ws.Range("D10").Value = 0.02030947775628627
$ws.Range("E10").Value = 0.001934235976789168
$ws.Range("F10").Value = 0.08123791102514506
$ws.Range("J10").Value = 0.1054158607350097
$ws.Range("O10").Value = 0.03288201160541586
$ws.Range("Q10").Value = 0.2263056092843327
$ws.Range("R10").Value = 0.08220502901353965
$ws.Range("S10").Value = 0.3646034816247582
$ws.Range("G11").Value = 0.1448598130841121
$ws.Range("J11").Value = 0.1261682242990654
$ws.Range("K11").Value = 0.2149532710280374
$ws.Range("L11").Value = 0.4953271028037383
$ws.Range("S11").Value = 0.01869158878504673
$ws.Range("G12").Value = 0.6635514018691588
$ws.Range("J12").Value = 0.2803738317757009
$ws.Range("K12").Value = 0.009345794392523364
$ws.Range("S12").Value = 0.04672897196261682
$ws.Range("G13").Value = 0.6976744186046512
$ws.Range("J13").Value = 0.2790697674418605
$ws.Range("S13").Value = 0.02325581395348837
$ws.Range("F15").Value = 0.00966183574879227
$ws.Range("H15").Value = 0.1835748792270532
$ws.Range("I15").Value = 0.04830917874396135
$ws.Range("J15").Value = 0.2705314009661836
$ws.Range("K15").Value = 0.03864734299516908
$ws.Range("M15").Value = 0.01932367149758454
$ws.Range("O15").Value = 0.03381642512077294
$ws.Range("S15").Value = 0.3961352657004831
$ws.Range("F16").Value = 0.006622516556291391
$ws.Range("H16").Value = 0.1788079470198675
$ws.Range("I16").Value = 0.1059602649006623
$ws.Range("J16").Value = 0.4039735099337748
$ws.Range("K16").Value = 0.1258278145695364
$ws.Range("M16").Value = 0.01324503311258278
$ws.Range("O16").Value = 0.05960264900662252
$ws.Range("S16").Value = 0.1059602649006623
$ws.Range("F17").Value = 0.01225490196078431
$ws.Range("H17").Value = 0.1764705882352941
$ws.Range("I17").Value = 0.1519607843137255
$ws.Range("J17").Value = 0.3995098039215687
$ws.Range("K17").Value = 0.05882352941176471
$ws.Range("M17").Value = 0.02450980392156863
$ws.Range("O17").Value = 0.06127450980392157
$ws.Range("S17").Value = 0.1151960784313725
$ws.Range("F18").Value = 0.02173913043478261
$ws.Range("H18").Value = 0.2028985507246377
$ws.Range("I18").Value = 0.1594202898550725
$ws.Range("J18").Value = 0.3768115942028986
$ws.Range("K18").Value = 0.05797101449275362
$ws.Range("M18").Value = 0.01449275362318841
$ws.Range("O18").Value = 0.07246376811594203
$ws.Range("S18").Value = 0.09420289855072464
$ws.Range("F19").Value = 0.01508429458740018
$ws.Range("H19").Value = 0.1925465838509317
$ws.Range("I19").Value = 0.1224489795918367
$ws.Range("J19").Value = 0.3815439219165927
$ws.Range("K19").Value = 0.09228039041703638
$ws.Range("M19").Value = 0.02218278615794144
$ws.Range("N19").Value = 0.0008873114463176575
$ws.Range("O19").Value = 0.06122448979591837
$ws.Range("S19").Value = 0.1118012422360248
